$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the text format ("@", numFmtId 49) on the header / row-label cells.
# These cells already carry this format; touching it again mirrors the
# original author re-writing the sheet (same effective style, new xf slot).
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"

# The substantive data change: the prediction/score value in B2.
$ws.Range("B2").Value = 3.2134672807675386
